$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update JD between reaction sets (B), metabolite sets (C), and gene sets (D)
# for draft network rows 4-32 with revised data.
$ws.Range("B4").Value = 0.73867975995635604
$ws.Range("C4").Value = 0.63328075709779197
$ws.Range("D4").Value = 0.52376137512638998
$ws.Range("B5").Value = 0.74064171122994704
$ws.Range("C5").Value = 0.63706563706563701
$ws.Range("D5").Value = 0.50712830957230104
$ws.Range("B6").Value = 0.59576968272620501
$ws.Range("C6").Value = 0.38789546079779902
$ws.Range("D6").Value = 0.55936675461741403
$ws.Range("B7").Value = 0.66258741258741305
$ws.Range("C7").Value = 0.53547776726584695
$ws.Range("D7").Value = 0.47149643705463201
$ws.Range("B8").Value = 0.74074074074074103
$ws.Range("C8").Value = 0.63835810332625598
$ws.Range("D8").Value = 0.48780487804878098
$ws.Range("B9").Value = 0.55791962174940901
$ws.Range("C9").Value = 0.37260273972602698
$ws.Range("D9").Value = 0.46245059288537599
$ws.Range("B10").Value = 0.65585851142225504
$ws.Range("C10").Value = 0.55301102629346899
$ws.Range("D10").Value = 0.37679558011049702
$ws.Range("B11").Value = 0.78462874511506697
$ws.Range("C11").Value = 0.71920932589964504
$ws.Range("D11").Value = 0.429973238180196
$ws.Range("B12").Value = 0.80955777460770295
$ws.Range("C12").Value = 0.69844961240310099
$ws.Range("D12").Value = 0.577596266044341
$ws.Range("B13").Value = 0.831398900427611
$ws.Range("C13").Value = 0.73723118279569899
$ws.Range("D13").Value = 0.577596266044341
$ws.Range("B14").Value = 0.81136543014996099
$ws.Range("C14").Value = 0.69211822660098499
$ws.Range("D14").Value = 0.577596266044341
$ws.Range("B15").Value = 0.81279999999999997
$ws.Range("C15").Value = 0.69256198347107401
$ws.Range("D15").Value = 0.577596266044341
$ws.Range("B16").Value = 0.80969351389878796
$ws.Range("C16").Value = 0.69844961240310099
$ws.Range("D16").Value = 0.577596266044341
$ws.Range("B17").Value = 0.80188124632569102
$ws.Range("C17").Value = 0.70647098065376901
$ws.Range("D17").Value = 0.47239263803680998
$ws.Range("B18").Value = 0.79814814814814805
$ws.Range("C18").Value = 0.70456092579986396
$ws.Range("D18").Value = 0.47239263803680998
$ws.Range("B19").Value = 0.87261580381471404
$ws.Range("C19").Value = 0.65672844480257897
$ws.Range("D19").Value = 0.49092970521541901
$ws.Range("B20").Value = 0.87195121951219501
$ws.Range("C20").Value = 0.65592264302981496
$ws.Range("D20").Value = 0.49322799097065501
$ws.Range("B21").Value = 0.86834733893557403
$ws.Range("C21").Value = 0.64516129032258096
$ws.Range("D21").Value = 0.5
$ws.Range("B22").Value = 0.86883116883116895
$ws.Range("C22").Value = 0.65919629057187001
$ws.Range("D22").Value = 0.46502057613168701
$ws.Range("B23").Value = 0.86860841423948199
$ws.Range("C23").Value = 0.65426356589147305
$ws.Range("D23").Value = 0.46240988671472699
$ws.Range("B24").Value = 0.86481113320079495
$ws.Range("C24").Value = 0.64909520062942605
$ws.Range("D24").Value = 0.46632653061224499
$ws.Range("B25").Value = 0.83909895414320201
$ws.Range("C25").Value = 0.72174590802805905
$ws.Range("D25").Value = 0.52835051546391798
$ws.Range("B26").Value = 0.83845547675334897
$ws.Range("C26").Value = 0.72499999999999998
$ws.Range("D26").Value = 0.53366583541147095
$ws.Range("B27").Value = 0.83863080684596603
$ws.Range("C27").Value = 0.71927042030134802
$ws.Range("D27").Value = 0.53174603174603197
$ws.Range("B28").Value = 0.83609576427256005
$ws.Range("C28").Value = 0.69216417910447803
$ws.Range("D28").Value = 0.59777102330293796
$ws.Range("B29").Value = 0.83827061649319501
$ws.Range("C29").Value = 0.72006220839813395
$ws.Range("D29").Value = 0.52600170502983801
$ws.Range("B30").Value = 0.84727551803530299
$ws.Range("C30").Value = 0.57965451055662198
$ws.Range("D30").Value = 0.405275779376499
$ws.Range("B31").Value = 0.879525593008739
$ws.Range("C31").Value = 0.71601615074024205
$ws.Range("D31").Value = 0.57185628742515004
$ws.Range("B32").Value = 0.87756370416407703
$ws.Range("C32").Value = 0.72015915119363405
$ws.Range("D32").Value = 0.56539235412474897

# Update view state: scroll position and active selection
$ws.Range("C12").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save() | Out-Null
